$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 10.97115645772287
$ws.Range("C2").Value = 11.562596071250281
$ws.Range("D2").Value = 8.5555184204074042
$ws.Range("E2").Value = 0.82064937990615228

$ws.Range("B3").Value = 28.480775705120902
$ws.Range("C3").Value = 4.1569834129557313
$ws.Range("D3").Value = 2.3063268374120849
$ws.Range("E3").Value = 2.532526546573743

$ws.Range("B1:E3").Select()
